$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing email list. This shifts the
# existing rows (A1:A5) down to (A2:A6) while keeping their original
# content/shared-string values intact.
$ws.Rows(1).Insert()

# The newly inserted row 1 has no formatting yet; copy the formatting that
# row 2 already carries (the original row 1's style) onto the new row 1 so
# the header cell matches the rest of the column's style.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Finally, give the new header cell its text.
$ws.Range("A1").Value = "Email"
